# Penalty Reward System (unfinished) — per commit message.
# 1) On the "Weekly Quantity" sheet, remove the obsolete first data row
#    (week 45368.99999999999 / qty 20) so every following row shifts up.
# 2) On the "Monthly Trend" sheet, correct the requested quantity for the
#    45382.99999999999 month from 250 to 230.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows("2:2").Delete()

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B2").Value = 230
